# Apply per-row cryptocurrency price/volume updates (GitHub Actions scrape refresh).
# D (Price) values that look numeric are written with a leading apostrophe so
# Excel keeps them as literal text (matching the workbook's inlineStr storage)
# instead of silently coercing them to the Number type.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.579.93"
$ws.Range("E2").Value = "  -0.73%  "

$ws.Range("D3").Value = "3.770.15"
$ws.Range("E3").Value = "  -1.91%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("D5").Value = "'596.18"
$ws.Range("E5").Value = "  -0.76%  "

$ws.Range("D6").Value = "'170.25"
$ws.Range("E6").Value = "  +1.49%  "

$ws.Range("D7").Value = "3.768.19"
$ws.Range("E7").Value = "  -1.90%  "

$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("D9").Value = "'0.528"
$ws.Range("E9").Value = "  -0.11%  "

$ws.Range("E10").Value = "  +1.01%  "

$ws.Range("D11").Value = "'6.48"
$ws.Range("E11").Value = "  +0.25%  "

$ws.Range("D12").Value = "'0.456"
$ws.Range("E12").Value = "  -0.53%  "

$ws.Range("D13").Value = "'0.0000278"
$ws.Range("E13").Value = "  +6.87%  "

$ws.Range("D14").Value = "'36.76"
$ws.Range("E14").Value = "  -0.84%  "

$ws.Range("D15").Value = "4.403.76"
$ws.Range("E15").Value = "  -1.86%  "

$ws.Range("D16").Value = "3.765.11"
$ws.Range("E16").Value = "  -1.99%  "

$ws.Range("D17").Value = "'18.89"
$ws.Range("E17").Value = "  +3.07%  "

$ws.Range("D18").Value = "67.606.76"
$ws.Range("E18").Value = "  -0.78%  "

$ws.Range("D19").Value = "'7.25"
$ws.Range("E19").Value = "  -2.03%  "

$ws.Range("E20").Value = "  +0.88%  "

$ws.Range("D21").Value = "'10.62"
$ws.Range("E21").Value = "  -4.23%  "

$ws.Range("D22").Value = "'469.87"
$ws.Range("E22").Value = "  +0.67%  "

$ws.Range("E23").Value = "  -1.65%  "

$ws.Range("D24").Value = "'0.0000148"
$ws.Range("E24").Value = "  -7.33%  "

$ws.Range("E25").Value = "  +1.29%  "

$ws.Range("E26").Value = "  +0.54%  "

$ws.Range("E28").Value = "  +3.70%  "

$ws.Range("E29").Value = "  +0.13%  "

$ws.Range("D30").Value = "'2.91"
$ws.Range("E30").Value = "  -1.89%  "

$ws.Range("D31").Value = "3.913.74"
$ws.Range("E31").Value = "  -1.94%  "

$ws.Range("D32").Value = "'7.73"
$ws.Range("E32").Value = "  +1.00%  "

$ws.Range("D33").Value = "'2.25"
$ws.Range("E33").Value = "  -2.85%  "

$ws.Range("D34").Value = "'30.52"
$ws.Range("E34").Value = "  -2.52%  "

$ws.Range("E35").Value = "  -4.08%  "

$ws.Range("D36").Value = "3.735.91"
$ws.Range("E36").Value = "  -1.80%  "

$ws.Range("D37").Value = "'3.86"
$ws.Range("E37").Value = "  +6.27%  "

$ws.Range("E38").Value = "  +1.32%  "

$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").Value = "'0.138"
$ws.Range("E39").Value = "  -1.52%  "

$ws.Range("B40").Value = "Filecoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D40").Value = "'5.90"
$ws.Range("E40").Value = "  -0.49%  "

$ws.Range("E41").Value = "  -1.31%  "

$ws.Range("D42").Value = "'1.00"
$ws.Range("E42").Value = "  +0.09%  "

$ws.Range("E43").Value = "  +0.37%  "

$ws.Range("E45").Value = "  +0.67%  "

$ws.Range("E46").Value = "  -1.03%  "

$ws.Range("D47").Value = "'45.91"
$ws.Range("E47").Value = "  -2.24%  "

$ws.Range("D48").Value = "'399.15"
$ws.Range("E48").Value = "  -5.36%  "

$ws.Range("D49").Value = "'0.000270"
$ws.Range("E49").Value = "  -7.15%  "

$ws.Range("D50").Value = "'141.91"
$ws.Range("E50").Value = "  -0.32%  "

$ws.Range("E51").Value = "  -0.31%  "
